$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15

$ws.Cells.Item($row, 1).Value = "Record"
$ws.Cells.Item($row, 2).Value = "Balanço Geral"
$ws.Cells.Item($row, 3).Value = "Trânsito"

# Ensure the date-like text in column D is stored as plain text rather than
# being auto-converted to a date/number by Excel.
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "2025-04-01T11:58"

$ws.Cells.Item($row, 5).Value = "Neutro"
$ws.Cells.Item($row, 6).Value = "Carreta invade o calçadão, derruba poste e deixa o Centro sem energia. Repórter *ao vivo*. Previsão é o caminhão ser retirado até 6h da tarde. Local isolado. Equipe da Enel no local. Guarda Municipal fez isolamento da área central para evitar déficit. Motorista é de Caxias do Sul (SC) e estava indo levar mercadoria para São Pedro da Aldeia. GPS estava marcando que aqui era uma rua. Estava escuro no momento. Entrevista com motorista da carreta e com comerciantes que ficaram sem energia. "
